# Apply row-content permutation for rows 3-15 and 17-26 (row 16 unchanged)
# as described by the diff: species/observation records were reordered
# while columns C, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY stay in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 112379407
$ws.Range("B3").Value = 78713
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("Q3").Value = 616516
$ws.Range("R3").Value = 7220340

# Row 4
$ws.Range("A4").Value = 112379409
$ws.Range("B4").Value = 77053
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6437
$ws.Range("F4").Value = "Blanksvart spiklav"
$ws.Range("G4").Value = "Calicium denigratum"
$ws.Range("H4").Value = "(Vain.) Tibell"
$ws.Range("Q4").Value = 616503
$ws.Range("R4").Value = 7220346

# Row 5
$ws.Range("A5").Value = 112379420
$ws.Range("B5").Value = 78713
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("Q5").Value = 616535
$ws.Range("R5").Value = 7220125

# Row 6
$ws.Range("A6").Value = 112379391
$ws.Range("B6").Value = 89499
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 112
$ws.Range("F6").Value = "Stjärntagging"
$ws.Range("G6").Value = "Asterodon ferruginosus"
$ws.Range("H6").Value = "Pat."
$ws.Range("Q6").Value = 616775
$ws.Range("R6").Value = 7220246

# Row 7
$ws.Range("A7").Value = 112379404
$ws.Range("B7").Value = 77403
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 228912
$ws.Range("F7").Value = "Mörk kolflarnlav"
$ws.Range("G7").Value = "Carbonicola myrmecina"
$ws.Range("H7").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q7").Value = 616559
$ws.Range("R7").Value = 7220391

# Row 8
$ws.Range("A8").Value = 112379392
$ws.Range("B8").Value = 77650
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 616797
$ws.Range("R8").Value = 7220278

# Row 9
$ws.Range("A9").Value = 112379406
$ws.Range("B9").Value = 90113
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 760
$ws.Range("F9").Value = "Doftticka"
$ws.Range("G9").Value = "Haploporus odorus"
$ws.Range("H9").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q9").Value = 616530
$ws.Range("R9").Value = 7220354

# Row 10
$ws.Range("A10").Value = 112379399
$ws.Range("B10").Value = 90826
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 4366
$ws.Range("F10").Value = "Skarp dropptaggsvamp"
$ws.Range("G10").Value = "Hydnellum peckii"
$ws.Range("H10").Value = "Banker"
$ws.Range("Q10").Value = 616801
$ws.Range("R10").Value = 7220436

# Row 11
$ws.Range("A11").Value = 112379401
$ws.Range("B11").Value = 94301
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 53
$ws.Range("F11").Value = "Vedtrappmossa"
$ws.Range("G11").Value = "Crossocalyx hellerianus"
$ws.Range("H11").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q11").Value = 616593
$ws.Range("R11").Value = 7220379

# Row 12
$ws.Range("A12").Value = 112379396
$ws.Range("B12").Value = 79580
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 1049
$ws.Range("F12").Value = "Kortskaftad ärgspik"
$ws.Range("G12").Value = "Microcalicium ahlneri"
$ws.Range("H12").Value = "Tibell"
$ws.Range("Q12").Value = 616825
$ws.Range("R12").Value = 7220323

# Row 13
$ws.Range("A13").Value = 112379419
$ws.Range("B13").Value = 78713
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6458
$ws.Range("F13").Value = "Lunglav"
$ws.Range("G13").Value = "Lobaria pulmonaria"
$ws.Range("H13").Value = "(L.) Hoffm."
$ws.Range("Q13").Value = 616545
$ws.Range("R13").Value = 7220142

# Row 14
$ws.Range("A14").Value = 112379414
$ws.Range("B14").Value = 78713
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6458
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."
$ws.Range("Q14").Value = 616567
$ws.Range("R14").Value = 7220205

# Row 15
$ws.Range("A15").Value = 112379418
$ws.Range("B15").Value = 90113
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 760
$ws.Range("F15").Value = "Doftticka"
$ws.Range("G15").Value = "Haploporus odorus"
$ws.Range("H15").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q15").Value = 616544
$ws.Range("R15").Value = 7220144

# Row 17
$ws.Range("A17").Value = 112379410
$ws.Range("B17").Value = 78713
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 6458
$ws.Range("F17").Value = "Lunglav"
$ws.Range("G17").Value = "Lobaria pulmonaria"
$ws.Range("H17").Value = "(L.) Hoffm."
$ws.Range("Q17").Value = 616528
$ws.Range("R17").Value = 7220271

# Row 18
$ws.Range("A18").Value = 112379412
$ws.Range("B18").Value = 89553
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 1202
$ws.Range("F18").Value = "Ullticka"
$ws.Range("G18").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H18").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q18").Value = 616546
$ws.Range("R18").Value = 7220243

# Row 19
$ws.Range("A19").Value = 112379415
$ws.Range("B19").Value = 78713
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6458
$ws.Range("F19").Value = "Lunglav"
$ws.Range("G19").Value = "Lobaria pulmonaria"
$ws.Range("H19").Value = "(L.) Hoffm."
$ws.Range("Q19").Value = 616569
$ws.Range("R19").Value = 7220196

# Row 20
$ws.Range("A20").Value = 112379416
$ws.Range("B20").Value = 77650
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 616572
$ws.Range("R20").Value = 7220163

# Row 21
$ws.Range("A21").Value = 112379411
$ws.Range("B21").Value = 78713
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("Q21").Value = 616537
$ws.Range("R21").Value = 7220266

# Row 22
$ws.Range("A22").Value = 112379413
$ws.Range("B22").Value = 89553
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 1202
$ws.Range("F22").Value = "Ullticka"
$ws.Range("G22").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H22").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q22").Value = 616565
$ws.Range("R22").Value = 7220210

# Row 23
$ws.Range("A23").Value = 112379417
$ws.Range("B23").Value = 78746
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 6463
$ws.Range("F23").Value = "Bårdlav"
$ws.Range("G23").Value = "Nephroma parile"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("Q23").Value = 616541
$ws.Range("R23").Value = 7220145

# Row 24
$ws.Range("A24").Value = 112379397
$ws.Range("B24").Value = 77402
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 6446
$ws.Range("F24").Value = "Kolflarnlav"
$ws.Range("G24").Value = "Carbonicola anthracophila"
$ws.Range("H24").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q24").Value = 616824
$ws.Range("R24").Value = 7220321

# Row 25
$ws.Range("A25").Value = 112379393
$ws.Range("B25").Value = 89549
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 1108
$ws.Range("F25").Value = "Harticka"
$ws.Range("G25").Value = "Pelloporus leporinus"
$ws.Range("H25").Value = "(Fr.) Krieglst."
$ws.Range("Q25").Value = 616793
$ws.Range("R25").Value = 7220301

# Row 26
$ws.Range("A26").Value = 112379405
$ws.Range("B26").Value = 78242
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6453
$ws.Range("F26").Value = "Vedskivlav"
$ws.Range("G26").Value = "Hertelidea botryosa"
$ws.Range("H26").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q26").Value = 616558
$ws.Range("R26").Value = 7220391
